$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix label / comment text (underscores -> spaces, clarified wording) ---
$ws.Range("A3").Value = "Sore throat"
$ws.Range("A4").Value = "Muscle ache"
$ws.Range("A5").Value = "Lack of appetite"
$ws.Range("A6").Value = "Smell before illness"
$ws.Range("A7").Value = "Smell during illness"
$ws.Range("C6").Value = "Scale 1-10 (1- no sense of smell, 10 - excellent sence of smell)"
$ws.Range("C7").Value = "Scale 1-10 (1- no sense of smell, 10 - excellent sence of smell)"

# --- Fill in the questionnaire answers (D2:D7) ---
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 1

# --- D9 now reads the probability input from B12 instead of the derived odds in B11 ---
$ws.Range("D9").Formula = "=`$B`$12*EXP(D8)/(1+`$B`$12*EXP(D8))"

# --- D10 now renders the probability as a rounded percentage instead of the raw number ---
$ws.Range("D10").Formula = '=IF(D9<0.5,"Can not be determined",_xlfn.CONCAT(ROUND(D9*100,0), "%"))'

# --- E10 picked up a small red "note" style (new font + cell style) ---
$ws.Range("E10").Font.Size = 10
$ws.Range("E10").Font.Color = 255
$ws.Range("E10").HorizontalAlignment = -4131
$ws.Range("E10").VerticalAlignment = -4108

# --- Column C widened to fit the longer descriptive text ---
$ws.Columns.Item(3).ColumnWidth = 49.7

# --- Restore the selection to match the saved view ---
$ws.Range("D13").Select() | Out-Null

Write-Output "edit applied"
